$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 13
$ws.Cells.Item($row, 1).Value = "2025-08-11 17:47"
$ws.Cells.Item($row, 2).Value = "4becd59"
$ws.Cells.Item($row, 3).Value = "[FIX]: Resolve category arrow visibility issues - eliminate duplicate arrows in shared categories"
$ws.Cells.Item($row, 4).Value = 2
$ws.Cells.Item($row, 5).Value = 16
$ws.Cells.Item($row, 6).Value = 7
$ws.Cells.Item($row, 7).Value = "Category arrow visibility fix - eliminate duplicate arrows"
$ws.Cells.Item($row, 8).Value = "Local"
